$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The worksheet table has data in rows 1, 5, 9, 13, 17 (5 columns each = 25 cells),
# laid out in document order. Map each cell (by row/col) to its old -> new text.
$replacements = @(
  @{Row=1;  Col=1; Old="37÷8=4, 5";   New="12÷6=2, 0"},
  @{Row=1;  Col=2; Old="38÷2=19, 0";  New="60÷3=20, 0"},
  @{Row=1;  Col=3; Old="79÷2=39, 1";  New="51÷2=25, 1"},
  @{Row=1;  Col=4; Old="29÷2=14, 1";  New="99÷7=14, 1"},
  @{Row=1;  Col=5; Old="71÷3=23, 2";  New="61÷6=10, 1"},

  @{Row=5;  Col=1; Old="83÷8=10, 3";  New="26÷3=8, 2"},
  @{Row=5;  Col=2; Old="60÷8=7, 4";   New="33÷7=4, 5"},
  @{Row=5;  Col=3; Old="77÷4=19, 1";  New="20÷5=4, 0"},
  @{Row=5;  Col=4; Old="39÷5=7, 4";   New="29÷7=4, 1"},
  @{Row=5;  Col=5; Old="42÷6=7, 0";   New="21÷2=10, 1"},

  @{Row=9;  Col=1; Old="45÷9=5, 0";   New="59÷3=19, 2"},
  @{Row=9;  Col=2; Old="59÷4=14, 3";  New="82÷8=10, 2"},
  @{Row=9;  Col=3; Old="36÷9=4, 0";   New="21÷9=2, 3"},
  @{Row=9;  Col=4; Old="38÷2=19, 0";  New="33÷9=3, 6"},
  @{Row=9;  Col=5; Old="37÷3=12, 1";  New="70÷9=7, 7"},

  @{Row=13; Col=1; Old="47÷9=5, 2";   New="76÷5=15, 1"},
  @{Row=13; Col=2; Old="13÷5=2, 3";   New="33÷3=11, 0"},
  @{Row=13; Col=3; Old="38÷7=5, 3";   New="54÷9=6, 0"},
  @{Row=13; Col=4; Old="59÷4=14, 3";  New="22÷6=3, 4"},
  @{Row=13; Col=5; Old="19÷9=2, 1";   New="40÷2=20, 0"},

  @{Row=17; Col=1; Old="35÷8=4, 3";   New="23÷5=4, 3"},
  @{Row=17; Col=2; Old="21÷4=5, 1";   New="73÷4=18, 1"},
  @{Row=17; Col=3; Old="39÷4=9, 3";   New="48÷6=8, 0"},
  @{Row=17; Col=4; Old="10÷8=1, 2";   New="72÷2=36, 0"},
  @{Row=17; Col=5; Old="62÷2=31, 0";  New="36÷3=12, 0"}
)

foreach ($item in $replacements) {
  $cell = $t.Cell($item.Row, $item.Col)
  # Scope the Find/Replace to exactly this cell's range (recomputed fresh each
  # time so earlier replacements that change text length don't throw off
  # later cell boundaries), and replace only the first (only) match inside it.
  $rng = $d.Range($cell.Range.Start, $cell.Range.End)
  $rng.Find.Execute($item.Old, $true, $false, $false, $false, $false, $true, 0, $false, $item.New, 1) | Out-Null
}
